# Práctica 3 - Product Backlog.xlsx
# "Elaboración de plantilla para sobre 1 y dossier publicitario"
#
# Fill in ESTIMACIÓN (column D) and SPRINT (column F) values for the
# backlog items, add a totals row with a SUM formula, resize row 1,
# and update the view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (ESTIMACIÓN, min) updates -----------------------------------
$ws.Range("D2").Value  = 35
$ws.Range("D3").Value  = 5
$ws.Range("D4").Value  = 10
$ws.Range("D5").Value  = 15
$ws.Range("D6").Value  = 15
$ws.Range("D7").Value  = 25
$ws.Range("D8").Value  = 10
$ws.Range("D9").Value  = 10
$ws.Range("D10").Value = 10
$ws.Range("D11").Value = 15
$ws.Range("D12").Value = 15
$ws.Range("D13").Value = 15
$ws.Range("D14").Value = 5
$ws.Range("D15").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("D21").Value = 7
$ws.Range("D22").Value = 20

# --- Column F (SPRINT) updates ---------------------------------------------
$ws.Range("F2").Value  = 1
$ws.Range("F3").Value  = 1
$ws.Range("F4").Value  = 1
$ws.Range("F5").Value  = 1
$ws.Range("F6").Value  = 1
$ws.Range("F7").Value  = 2
$ws.Range("F8").Value  = 2
$ws.Range("F9").Value  = 2
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 3
$ws.Range("F17").Value = 3
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 3
$ws.Range("F20").Value = 3
$ws.Range("F21").Value = 3
$ws.Range("F22").Value = 3

# --- New totals row ----------------------------------------------------
$ws.Range("D23").Formula = "=SUM(D2:D22)"

# --- Row height for header row ------------------------------------------
$ws.Rows.Item(1).RowHeight = 28.2

# --- View changes: zoom + selection, drop frozen/scrolled top-left -----
$excel.ActiveWindow.Zoom = 91
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select()
